$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row to reflect the new column layout / names.
# The order of assignment below matches the order in which the shared
# strings were (re)written by the original author, so that newly
# introduced strings land at the same shared-string-table positions.
$ws.Range("A1").Value = "name"
$ws.Range("C1").Value = "state"
$ws.Range("B1").Value = "city"
$ws.Range("E1").Value = "institute_type_id"
$ws.Range("H1").Value = "qs_rank"
$ws.Range("L1").Value = "shortnote"
$ws.Range("D1").Value = "views"
$ws.Range("F1").Value = "rating"
$ws.Range("G1").Value = "rank"
$ws.Range("I1").Value = "times_rank"
$ws.Range("J1").Value = "latitude_longitude"
$ws.Range("K1").Value = "featured"
$ws.Range("M1").Value = "established_year"

# Apply an explicit black font color to the last header cell
$ws.Range("M1").Font.Color = 0

# Set page orientation to portrait for printing
$ws.PageSetup.Orientation = 1
